$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the title text in A1 (backed by the shared-string table) ---
$ws.Range("A1").Value = "8.10.2.2 Камсыздандыруу компаниялардын финансылык көрсөткүчтөрү"

# xlPasteFormats: copy only the formatting (number format, font, borders,
# alignment, ...) of a cell, so the style index is reused instead of a new
# one being created.
$xlPasteFormats = -4122

# --- Add a new "2023" column (Q), mirroring the formatting of column P ---
# Row 3: year header
$ws.Range("P3").Copy() | Out-Null
$ws.Range("Q3").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("Q3").Value = 2023

# Row 4: number of reporting insurance companies
$ws.Range("P4").Copy() | Out-Null
$ws.Range("Q4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("Q4").Value = 16

# Row 5: insurance premiums (mln. soms)
$ws.Range("P5").Copy() | Out-Null
$ws.Range("Q5").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("Q5").Value = 3031.4

$excel.CutCopyMode = 0

# --- Leave the selection on the default top-left cell ---
$ws.Range("A1").Select() | Out-Null
